$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "4_201115_0909_rf_with_3in1_added_profitloss_relatd_features"
$ws.Range("B5").Value = "0.57+"
$ws.Range("C5").Value = "random foreset, in modifed 3in1 data set (new profit loss related features), train on 1-34, valid on 41 - 43 (.88)"

$ws.Range("C6").Select()
